$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Mobile number" column (old column M); shifts N:U left by one.
$ws.Columns("M").Delete()

# --- Row 4 header restyle ---
# A4:I4 / J4:T4 used an italic-only style; bold it to match the new header look.
$ws.Range("A4:I4").Font.Bold = $true
$ws.Range("J4:T4").Font.Bold = $true

# New "Time Off Reason" header in column S (old T, now S after the column delete).
$ws.Range("S4").Value = "Time Off Reason"

# --- Row 5 new data (second staff record) ---
$ws.Range("J5").Value = "Monnyka Pin"

$ws.Range("K5:M5").NumberFormat = "@"
$ws.Range("K5").Value = "QA Tester"
$ws.Range("L5").Value = "Male"
$ws.Range("M5").Value = "Cambodia"

$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "01/08/2017"

$ws.Range("P5").Value = "Oliver"
$ws.Range("Q5").Value = "Phnom Penh"

$ws.Range("R5:S5").NumberFormat = "@"
$ws.Range("R5").Value = "Thank for your hard working on the stars app feature, especially try to finish all task on time."
$ws.Range("S5").Value = "I have to go to the bank tomorrow "

# Selection ends up parked on S9 after the edits.
$ws.Range("S9").Select()
